# feat: add 2022-Q1 data
#
# 1) Insert a brand-new "2022-Q1" worksheet right before the "总计" sheet,
#    populated with the per-fund holding breakdown for 2022-Q1.
# 2) Update the "总计" (totals) sheet by inserting a new row for the
#    2022-Q1 quarter at the top of its data block and renumbering the
#    existing rows' running index.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q1" sheet
# ---------------------------------------------------------------------

$totals = $wb.Worksheets.Item("总计")
$src = $wb.Worksheets.Item("2021-Q4")   # donor sheet for formatting (same layout)

$q1 = $wb.Worksheets.Add($totals)
$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Copy the header formatting (bold + border, style index 2) from the
# donor sheet so the new header matches the rest of the workbook.
$src.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$rows = @(
    @{ code = "005014"; name = "泰康景泰回报混合A";                 size = "11.64"; pos = "27.21"; ratio = "1.39"; value = "0.1618"; rank = 7  },
    @{ code = "005775"; name = "中加转型动力灵活配置混合A";         size = "3.41";  pos = "66.34"; ratio = "2.05"; value = "0.0699"; rank = 10 },
    @{ code = "009927"; name = "工银瑞信聚利18个月定期开放混合A";   size = "5.54";  pos = "23.27"; ratio = "0.89"; value = "0.0493"; rank = 5  },
    @{ code = "005562"; name = "创金合信中证红利低波动指数C";       size = "1.85";  pos = "94.46"; ratio = "2.30"; value = "0.0426"; rank = 7  },
    @{ code = "005776"; name = "中加转型动力灵活配置混合C";         size = "1.92";  pos = "66.34"; ratio = "2.05"; value = "0.0394"; rank = 10 },
    @{ code = "512890"; name = "华泰柏瑞中证红利低波动ETF";         size = "1.36";  pos = "99.24"; ratio = "2.43"; value = "0.0330"; rank = 7  },
    @{ code = "005561"; name = "创金合信中证红利低波动指数A";       size = "1.22";  pos = "94.46"; ratio = "2.30"; value = "0.0281"; rank = 7  },
    @{ code = "005015"; name = "泰康景泰回报混合C";                 size = "0.63";  pos = "27.21"; ratio = "1.39"; value = "0.0088"; rank = 7  },
    @{ code = "009928"; name = "工银瑞信聚利18个月定期开放混合C";   size = "0.83";  pos = "23.27"; ratio = "0.89"; value = "0.0074"; rank = 5  }
)

$r = 2
foreach ($item in $rows) {
    $q1.Cells.Item($r, 1).Value = $r - 2
    $q1.Cells.Item($r, 2).Value = "'" + $item.code
    $q1.Cells.Item($r, 3).Value = $item.name
    $q1.Cells.Item($r, 4).Value = "'" + $item.size
    $q1.Cells.Item($r, 5).Value = "'" + $item.pos
    $q1.Cells.Item($r, 6).Value = "'" + $item.ratio
    $q1.Cells.Item($r, 7).Value = "'" + $item.value
    $q1.Cells.Item($r, 8).Value = $item.rank
    $r = $r + 1
}

# The apostrophe-prefix trick above forces text storage for
# numeric-looking strings (fund code / size / position / ratio /
# value), but it also stamps a transient "quote prefix" style on those
# cells. Strip that back off so the cells keep the workbook's default
# (no explicit style), matching the other quarter sheets.
$q1.Range("B2:G10").ClearFormats()

# Column A (running index) uses the same bold/border style as the
# other per-fund sheets.
$src.Range("A2").Copy()
$q1.Range("A2:A10").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Part 2: update "总计" with the new 2022-Q1 row
# ---------------------------------------------------------------------

# NOTE: re-resolve "总计" by name rather than reusing the handle
# captured above — inserting/renaming the new "2022-Q1" sheet shifts
# worksheet positions, and the old handle now resolves to "2022-Q1"
# instead of "总计".
$totals = $wb.Worksheets.Item("总计")

$totals.Rows.Item(2).Insert()
$totals.Rows.Item(2).ClearFormats()

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 9
$totals.Range("D2").Value = 0.44

$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)

for ($i = 3; $i -le 7; $i++) {
    $totals.Cells.Item($i, 1).Value = $i - 2
}

Write-Output "done"
